# Added New method for SelectDropdown using index
#
# tc016 / tc017 / tc018 previously drove their dropdown controls by setting
# the shared-string *label* of the option (e.g. "Bug", "Requirement - Wrong",
# "Mohit Aman", ...) into C2:G2. The new SelectDropdown-by-index method
# instead writes the numeric index of the option, so those five cells become
# plain numbers. This also removes "Requirement - Wrong" / "Mohit Aman" as
# needed strings; Excel automatically drops unused shared-string entries and
# renumbers the rest on save.

$wb = $excel.ActiveWorkbook

$sheetSelections = @{
    "tc016" = "E12"
    "tc017" = "F7"
    "tc018" = "E9"
}

foreach ($sheetName in @("tc016", "tc017", "tc018")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Activate()

    $ws.Range("C2").Value = 1
    $ws.Range("D2").Value = 1
    $ws.Range("E2").Value = 2
    $ws.Range("F2").Value = 1
    $ws.Range("G2").Value = 12

    [void]$ws.Range($sheetSelections[$sheetName]).Select()
}
